# dadosenade.xlsx - "Adiciona análise ENADE com gráficos e dados do Excel"
#
# Fills in the previously-incomplete 2026 projection row for the UFC /
# Fortaleza (ECOMP) cohort, completes the "CONCEITO" values that were
# missing for the two other 2026 rows in that block of the sheet, marks
# an empty trend cell with the underlined style used elsewhere, rebuilds
# the final SLOPE/INTERCEPT/FORECAST formulas against the full data range,
# nudges the active view/selection, and sets up the page for printing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 40 (UFC / ECOMP / Presencial / Fortaleza, year 2026) ---
# CONCEITO was missing, and the trend columns (INCLINAÇÃO/INTERSECÇÃO/
# PREVISÃO/CONCEITO FAIXA) had not yet been filled in for this cohort.
$ws.Range("B40").Value = 3
$ws.Range("G40").Value = 0.17857139999999999
$ws.Range("H40").Value = -357.32142900000002
$ws.Range("I40").Value = 4.4642860000000004
$ws.Range("J40").Value = 3

# --- Row 50 (UNIFOR / ECOMP / Presencial / Fortaleza, year 2026) ---
# Only the CONCEITO value was missing; the trend columns were already set.
$ws.Range("B50").Value = 4

# --- Row 58 (UFRN / ECOMP / Presencial / Natal, year 2026) ---
# Only the CONCEITO value was missing; the trend columns were already set.
$ws.Range("B58").Value = 3

# --- Row 63 (UNP / ECOMP / Presencial / Natal, year 2017) ---
# Mark the empty INCLINAÇÃO cell with the underlined font style.
$ws.Range("G63").Font.Underline = $true

# --- Row 66 totals: recompute against the full A2:A66 / B2:B66 range ---
$ws.Range("G66").Formula = "=SLOPE(B2:B66,A2:A66)"
$ws.Range("H66").Formula = "=INTERCEPT(B2:B66,A2:A66)"
$ws.Range("I66").Formula = "=FORECAST(2026,B2:B66,A2:A66)"

# --- View / selection adjustments ---
$ws.Range("A46").Select()
$excel.ActiveWindow.ScrollRow = 46
$ws.Range("A66").Select()

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
